# Append a new scrape batch (2025-10-16 12:39:02 JST) to the "ランサーズ" sheet.
# - refresh the "取得日時" timestamp on every existing data row (2-12)
# - fix a title/URL mix-up between rows 6 and 7 (their text got swapped)
# - append a brand new row 13 with its own hyperlink

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-10-16 12:39:02"

# Refresh the retrieval timestamp for every existing row.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Rows 6 and 7 had their title (B) and URL (F) text swapped in this edit.
$ws.Range("B6").Value = "大手クレジットカード企業向け、Google Cloudを利用したアジャイル開発共通基盤案件"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5414353"

$ws.Range("B7").Value = "大手クレジットカード企業向け、Google Cloudを利用したアジャイル開発共通基盤案件_ワーカー"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5414354"

# Append the new row 13.
$ws.Range("A13").Value = $newTimestamp
$ws.Range("B13").Value = "PowerAutomateメール監視して件名と本文内の条件一致時、社内システム操作&メール転送したい"
$ws.Range("C13").Value = "システム開発"
$ws.Range("D13").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E13").Value = "期限情報なし"
$ws.Range("F13").Value = "https://www.lancers.jp/work/detail/5414579"
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5414579")
$ws.Range("F13").Style = "Hyperlink"
$ws.Range("G13").Value = 28

$wb.Save()
